$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Lũy kế tháng HỆ THỐNG")

# The "last_edited_time" column (D) shares one string across rows 6-13.
# Update it to the new timestamp for all rows that currently hold it.
for ($r = 6; $r -le 13; $r++) {
    $cell = $ws.Range("D$r")
    if ($cell.Text -eq "2024-08-24T20:33:00.000Z") {
        $cell.Value = "2024-08-26T17:26:00.000Z"
    }
}

# Update the numeric figures on row 6.
$ws.Range("T6").Value = 118000000
$ws.Range("W6").Value = 198433000
$ws.Range("AA6").Value = 303617000
$ws.Range("AE6").Value = 502050000
$ws.Range("AH6").Value = 413050000
$ws.Range("AK6").Value = 66
$ws.Range("AN6").Value = 89000000
$ws.Range("AQ6").Value = 531050000
